{"js": "// Update master to output generated at 4250d90\n// Replace the worksheet date and each 3-digit x 1-digit multiplication\n// problem's text with the next day's generated values.\nconst replacements = [\n  [\"2024-08-29 Thursday\", \"2024-08-30 Friday\"],\n  [\"370\u00d73=\", \"489\u00d74=\"],\n  [\"811\u00d73=\", \"394\u00d73=\"],\n  [\"800\u00d72=\", \"460\u00d74=\"],\n  [\"228\u00d73=\", \"848\u00d78=\"],\n  [\"717\u00d78=\", \"395\u00d74=\"],\n  [\"430\u00d73=\", \"289\u00d78=\"],\n  [\"822\u00d79=\", \"639\u00d76=\"],\n  [\"346\u00d75=\", \"208\u00d79=\"],\n  [\"935\u00d75=\", \"959\u00d72=\"],\n  [\"275\u00d72=\", \"721\u00d78=\"],\n  [\"830\u00d76=\", \"138\u00d72=\"],\n  [\"292\u00d78=\", \"659\u00d75=\"],\n  [\"248\u00d78=\", \"644\u00d78=\"],\n  [\"104\u00d75=\", \"468\u00d77=\"],\n  [\"801\u00d79=\", \"970\u00d72=\"],\n  [\"450\u00d72=\", \"126\u00d73=\"],\n  [\"345\u00d72=\", \"270\u00d79=\"],\n  [\"510\u00d74=\", \"564\u00d74=\"],\n  [\"300\u00d75=\", \"326\u00d72=\"],\n  [\"778\u00d76=\", \"195\u00d73=\"],\n  [\"624\u00d78=\", \"246\u00d72=\"],\n  [\"193\u00d78=\", \"304\u00d79=\"],\n  [\"218\u00d72=\", \"829\u00d78=\"],\n  [\"964\u00d73=\", \"414\u00d73=\"],\n  [\"671\u00d74=\", \"623\u00d74=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update master to output generated at 4250d90\n# Replaces the worksheet date and each 3-digit x 1-digit multiplication\n# problem with the next day's generated values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-08-29 Thursday\", \"2024-08-30 Friday\"),\n    @(\"370\u00d73=\", \"489\u00d74=\"),\n    @(\"811\u00d73=\", \"394\u00d73=\"),\n    @(\"800\u00d72=\", \"460\u00d74=\"),\n    @(\"228\u00d73=\", \"848\u00d78=\"),\n    @(\"717\u00d78=\", \"395\u00d74=\"),\n    @(\"430\u00d73=\", \"289\u00d78=\"),\n    @(\"822\u00d79=\", \"639\u00d76=\"),\n    @(\"346\u00d75=\", \"208\u00d79=\"),\n    @(\"935\u00d75=\", \"959\u00d72=\"),\n    @(\"275\u00d72=\", \"721\u00d78=\"),\n    @(\"830\u00d76=\", \"138\u00d72=\"),\n    @(\"292\u00d78=\", \"659\u00d75=\"),\n    @(\"248\u00d78=\", \"644\u00d78=\"),\n    @(\"104\u00d75=\", \"468\u00d77=\"),\n    @(\"801\u00d79=\", \"970\u00d72=\"),\n    @(\"450\u00d72=\", \"126\u00d73=\"),\n    @(\"345\u00d72=\", \"270\u00d79=\"),\n    @(\"510\u00d74=\", \"564\u00d74=\"),\n    @(\"300\u00d75=\", \"326\u00d72=\"),\n    @(\"778\u00d76=\", \"195\u00d73=\"),\n    @(\"624\u00d78=\", \"246\u00d72=\"),\n    @(\"193\u00d78=\", \"304\u00d79=\"),\n    @(\"218\u00d72=\", \"829\u00d78=\"),\n    @(\"964\u00d73=\", \"414\u00d73=\"),\n    @(\"671\u00d74=\", \"623\u00d74=\")\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n\nWrite-Output \"Replaced $($replacements.Count) text runs\"\n"}
